# Madrid lat/long workbook edit:
#  - delete the "Parque de las Avenidas" row (row 16)
#  - insert a new column C "Superficie km2" with area (km2) data for every barrio

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the "Parque de las Avenidas" row entirely (old row 16).
$ws.Rows(16).Delete()

# 2) Insert a new blank column before the old column C (Compose/lat/long block
#    shifts one column to the right: C->D, D->E, E->F).
$ws.Columns("C").Insert()

# 3) Header for the new column.
$ws.Range("C1").Value = "Superficie km2"

# 4) Area values (km2) for each barrio, in the same row order as the sheet
#    now that "Parque de las Avenidas" has been removed (rows 2-25).
$areas = @(
    1.471,
    1.032,
    0.592,
    0.742,
    0.947,
    0.445,
    0.75,
    0.64,
    1.025,
    0.49,
    1.9,
    0.643,
    0.87,
    0.771,
    0.852,
    1.598,
    0.52,
    0.773,
    1.708,
    1.43,
    0.762,
    1.707,
    1.788,
    2.16
)

for ($i = 0; $i -lt $areas.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $areas[$i]
}

$ws.Range("C2").Select()
